# Realestate Update resale numbers 2025-01-04 21:32
# Append a new data row (row 8) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$rowNum = 8

# Columns A (Date) and D (Week) look like a date / a plain number to Excel's
# auto-detection, so force them to be stored as text (matching the other
# rows) by applying a text number format before assigning the value, then
# resetting the style back to Normal so no stray style index is left on
# the cell (consistent with the unstyled cells used by the other rows).
$ws.Cells.Item($rowNum, 1).NumberFormat = "@"
$ws.Cells.Item($rowNum, 1).Value = "2025-01-04"
$ws.Cells.Item($rowNum, 1).Style = "Normal"

$ws.Cells.Item($rowNum, 2).Value = "21:32:42"
$ws.Cells.Item($rowNum, 3).Value = "Saturday"

$ws.Cells.Item($rowNum, 4).NumberFormat = "@"
$ws.Cells.Item($rowNum, 4).Value = "00"
$ws.Cells.Item($rowNum, 4).Style = "Normal"

$ws.Cells.Item($rowNum, 5).Value = 127730
$ws.Cells.Item($rowNum, 6).Value = 143633
$ws.Cells.Item($rowNum, 7).Value = 168459
$ws.Cells.Item($rowNum, 8).Value = 158358
$ws.Cells.Item($rowNum, 9).Value = -1
$ws.Cells.Item($rowNum, 10).Value = 142214
$ws.Cells.Item($rowNum, 11).Value = -1
$ws.Cells.Item($rowNum, 12).Value = -1
$ws.Cells.Item($rowNum, 13).Value = 192468
$ws.Cells.Item($rowNum, 14).Value = 114826
$ws.Cells.Item($rowNum, 15).Value = 45479
$ws.Cells.Item($rowNum, 16).Value = 28307
$ws.Cells.Item($rowNum, 17).Value = 63568
$ws.Cells.Item($rowNum, 18).Value = -1
$ws.Cells.Item($rowNum, 19).Value = 48380
$ws.Cells.Item($rowNum, 20).Value = -1
